$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1692.909
$ws.Range("J17").Value = 1692.909
$ws.Range("L17").Value = 5078.727000000001
$ws.Range("N17").Value = -5414.727000000001

$ws.Range("H33").Value = 66.666664
$ws.Range("I33").Value = 70
$ws.Range("K33").Value = 70
$ws.Range("M33").Value = 159

$ws.Range("H116").Value = 67001
$ws.Range("J116").Value = 119999.5
$ws.Range("L116").Value = 119999.5
$ws.Range("N116").Value = -126883.5

$ws.Range("H132").Value = 3962.1
$ws.Range("I132").Value = 3919.3333
$ws.Range("J132").Value = 4026.25
$ws.Range("K132").Value = 11757.9999
$ws.Range("L132").Value = 12078.75
$ws.Range("M132").Value = -9227.999899999999
$ws.Range("N132").Value = -17138.75

$ws.Range("H138").Value = 4020
$ws.Range("I138").Value = 2397.6667
$ws.Range("K138").Value = 7193.000100000001
$ws.Range("M138").Value = -2053.000100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5285.8
$ws.Range("I32").Value = 5285.8
$ws.Range("K32").Value = 5285.8
$ws.Range("M32").Value = -4998.8

$ws.Range("H74").Value = 4888.6523
$ws.Range("I74").Value = 5297.905
$ws.Range("K74").Value = 5297.905
$ws.Range("M74").Value = -4423.905

$ws.Range("H77").Value = 4888.6523
$ws.Range("I77").Value = 5297.905
$ws.Range("K77").Value = 26489.525
$ws.Range("M77").Value = -22121.525

$ws.Range("H132").Value = 3454.3635
$ws.Range("J132").Value = 3814
$ws.Range("L132").Value = 11442
$ws.Range("N132").Value = -16502

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 13104
$ws.Range("J20").Value = 24004
$ws.Range("L20").Value = 24004
$ws.Range("N20").Value = -24498

$ws.Range("H22").Value = 542.8570999999999
$ws.Range("I22").Value = 542.8570999999999
$ws.Range("K22").Value = 542.8570999999999
$ws.Range("M22").Value = -369.8570999999999

$ws.Range("H94").Value = 30000
$ws.Range("J94").Value = 30000
$ws.Range("L94").Value = 30000
$ws.Range("N94").Value = -30902

$ws.Range("H134").Value = 2623.75
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 54.18182
$ws.Range("I7").Value = 61.166668
$ws.Range("J7").Value = 45.8
$ws.Range("K7").Value = 61.166668
$ws.Range("L7").Value = 45.8
$ws.Range("M7").Value = 51.833332
$ws.Range("N7").Value = -271.8

$ws.Range("H58").Value = 3331.1365
$ws.Range("I58").Value = 3388.111
$ws.Range("J58").Value = 3074.75
$ws.Range("K58").Value = 3388.111
$ws.Range("L58").Value = 3074.75
$ws.Range("M58").Value = -3185.111
$ws.Range("N58").Value = -3480.75

$ws.Range("H86").Value = 17441.666
$ws.Range("I86").Value = 18623.125
$ws.Range("K86").Value = 18623.125
$ws.Range("M86").Value = -17500.125

$ws.Range("H89").Value = 17441.666
$ws.Range("I89").Value = 18623.125
$ws.Range("K89").Value = 93115.625
$ws.Range("M89").Value = -87499.625

$ws.Range("H94").Value = 1521
$ws.Range("I94").Value = 1911.3334
$ws.Range("J94").Value = 350
$ws.Range("K94").Value = 1911.3334
$ws.Range("L94").Value = 350
$ws.Range("M94").Value = -1460.3334
$ws.Range("N94").Value = -1252

$ws.Range("H109").Value = 50285
$ws.Range("J109").Value = 50285
$ws.Range("L109").Value = 50285
$ws.Range("N109").Value = -52365

$ws.Range("H132").Value = 5149.6
$ws.Range("I132").Value = 4750
$ws.Range("J132").Value = 5249.5
$ws.Range("K132").Value = 14250
$ws.Range("L132").Value = 15748.5
$ws.Range("M132").Value = -11720
$ws.Range("N132").Value = -20808.5

$ws.Range("H134").Value = 7112.4
$ws.Range("I134").Value = 8315.5
$ws.Range("J134").Value = 2300
$ws.Range("K134").Value = 24946.5
$ws.Range("L134").Value = 6900
$ws.Range("M134").Value = -22411.5
$ws.Range("N134").Value = -11970

$ws.Range("H136").Value = 3331.1365
$ws.Range("I136").Value = 3388.111
$ws.Range("J136").Value = 3074.75
$ws.Range("K136").Value = 10164.333
$ws.Range("L136").Value = 9224.25
$ws.Range("M136").Value = -7614.332999999999
$ws.Range("N136").Value = -14324.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 5269.6875
$ws.Range("J34").Value = 6541.75
$ws.Range("L34").Value = 19625.25
$ws.Range("N34").Value = -19793.25

$ws.Range("H39").Value = 16750
$ws.Range("J39").Value = 16750
$ws.Range("L39").Value = 50250
$ws.Range("N39").Value = -50838

$ws.Range("H55").Value = 13500
$ws.Range("J55").Value = 13500
$ws.Range("L55").Value = 40500
$ws.Range("N55").Value = -40854

$ws.Range("H56").Value = 10666.667
$ws.Range("I56").Value = 10666.667
$ws.Range("K56").Value = 10666.667
$ws.Range("M56").Value = -10136.667

$ws.Range("H82").Value = 10000
$ws.Range("J82").Value = 10000
$ws.Range("L82").Value = 30000
$ws.Range("N82").Value = -30812

$ws.Range("H85").Value = 10000
$ws.Range("J85").Value = 10000
$ws.Range("L85").Value = 30000
$ws.Range("N85").Value = -32808

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H45").Value = 25000
$ws.Range("J45").Value = 25000
$ws.Range("L45").Value = 25000
$ws.Range("N45").Value = -26118

$ws.Range("H80").Value = 10549.25
$ws.Range("I80").Value = 4873.5
$ws.Range("J80").Value = 16225
$ws.Range("K80").Value = 4873.5
$ws.Range("L80").Value = 16225
$ws.Range("M80").Value = -3875.5
$ws.Range("N80").Value = -18221

$ws.Range("H83").Value = 10549.25
$ws.Range("I83").Value = 4873.5
$ws.Range("J83").Value = 16225
$ws.Range("K83").Value = 24367.5
$ws.Range("L83").Value = 81125
$ws.Range("M83").Value = -19375.5
$ws.Range("N83").Value = -91109

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6660.577
$ws.Range("I132").Value = 1463.4375
$ws.Range("K132").Value = 4390.3125
$ws.Range("M132").Value = -1860.3125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3129.7144
$ws.Range("I132").Value = 955
$ws.Range("K132").Value = 2865
$ws.Range("M132").Value = -335

$ws.Range("H136").Value = 2863.4546
$ws.Range("I136").Value = 3222
$ws.Range("K136").Value = 9666
$ws.Range("M136").Value = -7116
